$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new Price (column D) / new Volume(1h) (column E) values.
# A $null Price means that row's Price cell is unchanged in this update.
$updates = @(
    @{ Row = 2;  D = "29.233.60";     E = "  -0.01%  " },
    @{ Row = 3;  D = "1.844.76";      E = "  +0.24%  " },
    @{ Row = 4;  D = "0.9993";        E = "  -0.02%  " },
    @{ Row = 5;  D = "242.81";        E = "  +0.85%  " },
    @{ Row = 6;  D = "0.6631";        E = "  -0.91%  " },
    @{ Row = 7;  D = $null;           E = "  -0.03%  " },
    @{ Row = 8;  D = "44.95";         E = "  +7.39%  " },
    @{ Row = 9;  D = "0.07442";       E = "  +0.36%  " },
    @{ Row = 10; D = "0.2956";        E = "  -0.18%  " },
    @{ Row = 11; D = "23.32";         E = "  +2.07%  " },
    @{ Row = 12; D = "0.07765";       E = "  +0.67%  " },
    @{ Row = 13; D = "1.846.07";      E = "  +8.06%  " },
    @{ Row = 14; D = "5.021";         E = "  -0.22%  " },
    @{ Row = 15; D = "0.6724";        E = "  -1.00%  " },
    @{ Row = 16; D = "83.50";         E = "  -3.15%  " },
    @{ Row = 17; D = $null;           E = "  -0.26%  " },
    @{ Row = 18; D = "0.000008733";   E = "  +5.94%  " },
    @{ Row = 19; D = "29.236.61";     E = "  +0.82%  " },
    @{ Row = 20; D = "2.099.14";      E = "  +2.18%  " },
    @{ Row = 21; D = "12.55";         E = "  +0.03%  " },
    @{ Row = 22; D = "226.99";        E = "  -0.88%  " },
    @{ Row = 23; D = $null;           E = "  +0.11%  " },
    @{ Row = 24; D = "7.179";         E = "  -1.03%  " },
    @{ Row = 25; D = $null;           E = "  -0.06%  " },
    @{ Row = 26; D = "158.62";        E = "  -0.97%  " },
    @{ Row = 27; D = "0.1407";        E = "  -0.55%  " },
    @{ Row = 28; D = "8.640";         E = "  -0.90%  " },
    @{ Row = 29; D = $null;           E = "  -0.03%  " },
    @{ Row = 30; D = "1.508";         E = "  +0.03%  " },
    @{ Row = 31; D = "4.139";         E = "  -1.65%  " },
    @{ Row = 32; D = "4.056";         E = "  -0.59%  " },
    @{ Row = 33; D = $null;           E = "  -0.70%  " },
    @{ Row = 34; D = "0.05332";       E = "  -0.53%  " },
    @{ Row = 35; D = $null;           E = "  +0.09%  " },
    @{ Row = 36; D = "0.7474";        E = "  -1.50%  " },
    @{ Row = 37; D = "1.157";         E = "  +1.83%  " },
    @{ Row = 38; D = "2.656";         E = "  -0.91%  " },
    @{ Row = 39; D = "1.314.49";      E = "  -1.36%  " },
    @{ Row = 40; D = "0.01803";       E = "  +0.04%  " },
    @{ Row = 41; D = "2.756";         E = "  +0.96%  " },
    @{ Row = 42; D = "6.402";         E = "  +6.73%  " },
    @{ Row = 43; D = "0.9041";        E = "  -1.90%  " },
    @{ Row = 44; D = "1.0000";        E = "  +0.01%  " },
    @{ Row = 45; D = "103.40";        E = "  -0.03%  " },
    @{ Row = 46; D = "1.995.47";      E = "  +2.91%  " },
    @{ Row = 47; D = "0.07920";       E = "  +0.09%  " },
    @{ Row = 48; D = "65.42";         E = "  +2.50%  " },
    @{ Row = 49; D = $null;           E = "  -1.85%  " },
    @{ Row = 51; D = "1.754";         E = "  -0.63%  " }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($r, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
        $cellD.Style = "Normal"
    }

    $cellE = $ws.Cells.Item($r, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $u.E
    $cellE.Style = "Normal"
}
